# Append 4 new order rows (HK1038-HK1041) to Sheet1, rows 55-58.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 55; A = "HK1038"; B = "12/17/2025"; C = "Pranav";       D = "Daily Veg subscription";      E = 8;  F = 1; G = 8;  H = "Accepted" },
    @{ Row = 56; A = "HK1039"; B = "12/17/2025"; C = "Shalini Raju"; D = "Veg Curries";                 E = 15; F = 2; G = 30; H = "Accepted" },
    @{ Row = 57; A = "HK1040"; B = "12/17/2025"; C = "Abhilasha";    D = "Dosa Batter with Chutney";    E = 10; F = 1; G = 10; H = "Accepted" },
    @{ Row = 58; A = "HK1041"; B = "12/17/2025"; C = "Pranav";       D = "Daily Veg Subscription";      E = 8;  F = 1; G = 8;  H = "Accepted" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.A

    # Column B holds a date formatted as text ("12/17/2025"), not a real
    # date serial, matching how the other rows in this sheet store it.
    # Force text via NumberFormat, assign, then clear the format again so
    # the cell is left with the default (unstyled) appearance.
    $bCell = $ws.Cells.Item($rowNum, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $r.B
    $bCell.ClearFormats()

    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
}
